$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.038545033472623
$ws.Range("D2").Value2 = 1.046711504663584
$ws.Range("E2").Value2 = 1.042197612342531
$ws.Range("F2").Value2 = 1.055272785028723
$ws.Range("I2").Value2 = 1.040260252184625
$ws.Range("J2").Value2 = 1.043641891145105
$ws.Range("K2").Value2 = 1.049476175024802
$ws.Range("L2").Value2 = 1.044974979350359
$ws.Range("M2").Value2 = 1.058013697628278
$ws.Range("N2").Value2 = 1.045123982337929

$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.039679757028968
$ws.Range("D3").Value2 = 1.047608311872301
$ws.Range("E3").Value2 = 1.043281222806352
$ws.Range("F3").Value2 = 1.056303129456566
$ws.Range("I3").Value2 = 1.040537057269649
$ws.Range("J3").Value2 = 1.044420523394394
$ws.Range("K3").Value2 = 1.05018445889748
$ws.Range("L3").Value2 = 1.045868657300915
$ws.Range("M3").Value2 = 1.058856897109628
$ws.Range("N3").Value2 = 1.045903720334323

$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.040413958874661
$ws.Range("D4").Value2 = 1.048188248346104
$ws.Range("E4").Value2 = 1.043982703540205
$ws.Range("F4").Value2 = 1.056969774036463
$ws.Range("I4").Value2 = 1.040714385149609
$ws.Range("J4").Value2 = 1.044923773137249
$ws.Range("K4").Value2 = 1.050641776001725
$ws.Range("L4").Value2 = 1.046446647016822
$ws.Range("M4").Value2 = 1.059401839922814
$ws.Range("N4").Value2 = 1.046407684749536

$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.040722608642757
$ws.Range("D5").Value2 = 1.048431967954657
$ws.Range("E5").Value2 = 1.044277681174757
$ws.Range("F5").Value2 = 1.05725001782769
$ws.Range("I5").Value2 = 1.040788506752245
$ws.Range("J5").Value2 = 1.04513520138069
$ws.Range("K5").Value2 = 1.05083379519859
$ws.Range("L5").Value2 = 1.046689567253874
$ws.Range("M5").Value2 = 1.059630775099488
$ws.Range("N5").Value2 = 1.046619413245333

$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.040774431814512
$ws.Range("D6").Value2 = 1.048472884523564
$ws.Range("E6").Value2 = 1.044327213634115
$ws.Range("F6").Value2 = 1.057297071199103
$ws.Range("I6").Value2 = 1.040800927051948
$ws.Range("J6").Value2 = 1.045170693031804
$ws.Range("K6").Value2 = 1.050866022176048
$ws.Range("L6").Value2 = 1.046730350731937
$ws.Range("M6").Value2 = 1.059669205014593
$ws.Range("N6").Value2 = 1.046654955298663

$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.040418083100945
$ws.Range("D7").Value2 = 1.048191505276189
$ws.Range("E7").Value2 = 1.043986644750116
$ws.Range("F7").Value2 = 1.056973518721013
$ws.Range("I7").Value2 = 1.040715377244145
$ws.Range("J7").Value2 = 1.04492659879349
$ws.Range("K7").Value2 = 1.050644342702462
$ws.Range("L7").Value2 = 1.046449893190381
$ws.Range("M7").Value2 = 1.059404899589652
$ws.Range("N7").Value2 = 1.046410514418532

$ws.Range("B8").Value2 = 1.019999999999999
$ws.Range("C8").Value2 = 1.038928527488995
$ws.Range("D8").Value2 = 1.047014658914671
$ws.Range("E8").Value2 = 1.04256375912986
$ws.Range("F8").Value2 = 1.055621006646839
$ws.Range("I8").Value2 = 1.040354169013944
$ws.Range("J8").Value2 = 1.043905153280715
$ws.Range("K8").Value2 = 1.049715747893336
$ws.Range("L8").Value2 = 1.045277060098
$ws.Range("M8").Value2 = 1.058298798308887
$ws.Range("N8").Value2 = 1.045387618335968

$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.036303390446635
$ws.Range("D9").Value2 = 1.044938166030716
$ws.Range("E9").Value2 = 1.040058827633832
$ws.Range("F9").Value2 = 1.053237264669981
$ws.Range("I9").Value2 = 1.039704016474935
$ws.Range("J9").Value2 = 1.042100800566902
$ws.Range("K9").Value2 = 1.048071859418064
$ws.Range("L9").Value2 = 1.043208221594899
$ws.Range("M9").Value2 = 1.056344618891654
$ws.Range("N9").Value2 = 1.04358070323434

$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.034553002282281
$ws.Range("D10").Value2 = 1.043551993294487
$ws.Range("E10").Value2 = 1.038390442821604
$ws.Range("F10").Value2 = 1.05164779480355
$ws.Range("I10").Value2 = 1.039261399945553
$ws.Range("J10").Value2 = 1.040894895729967
$ws.Range("K10").Value2 = 1.046970824719328
$ws.Range("L10").Value2 = 1.041827522310602
$ws.Range("M10").Value2 = 1.05503840171204
$ws.Range("N10").Value2 = 1.042373085874216

$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.033794979874562
$ws.Range("D11").Value2 = 1.042951326154262
$ws.Range("E11").Value2 = 1.037668379011511
$ws.Range("F11").Value2 = 1.050959460094206
$ws.Range("I11").Value2 = 1.039067563969267
$ws.Range("J11").Value2 = 1.040372007813037
$ws.Range("K11").Value2 = 1.046492849844693
$ws.Range("L11").Value2 = 1.041229308607544
$ws.Range("M11").Value2 = 1.054471978180415
$ws.Range("N11").Value2 = 1.041849455396469

$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.033513401295211
$ws.Range("D12").Value2 = 1.042728144484458
$ws.Range("E12").Value2 = 1.037400225404532
$ws.Range("F12").Value2 = 1.050703768920865
$ws.Range("I12").Value2 = 1.038995237000792
$ws.Range("J12").Value2 = 1.04017767473315
$ws.Range("K12").Value2 = 1.046315125102965
$ws.Range("L12").Value2 = 1.041007050410193
$ws.Range("M12").Value2 = 1.054261459228707
$ws.Range("N12").Value2 = 1.041654846341318

$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.033573801570936
$ws.Range("D13").Value2 = 1.04277602075318
$ws.Range("E13").Value2 = 1.037457742874793
$ws.Range("F13").Value2 = 1.050758616141385
$ws.Range("I13").Value2 = 1.039010766210174
$ws.Range("J13").Value2 = 1.040219364790502
$ws.Range("K13").Value2 = 1.046353255977149
$ws.Range("L13").Value2 = 1.041054728046782
$ws.Range("M13").Value2 = 1.054306621877743
$ws.Range("N13").Value2 = 1.041696595603334

$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.033771704821509
$ws.Range("D14").Value2 = 1.042932879259184
$ws.Range("E14").Value2 = 1.037646212272828
$ws.Range("F14").Value2 = 1.050938324850044
$ws.Range("I14").Value2 = 1.039061592084294
$ws.Range("J14").Value2 = 1.040355946411568
$ws.Range("K14").Value2 = 1.046478162809428
$ws.Range("L14").Value2 = 1.041210937801141
$ws.Range("M14").Value2 = 1.054454579150774
$ws.Range("N14").Value2 = 1.041833371185967

$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.03389363747668
$ws.Range("D15").Value2 = 1.043029516055809
$ws.Range("E15").Value2 = 1.037762341496298
$ws.Range("F15").Value2 = 1.051049047579144
$ws.Range("I15").Value2 = 1.039092864161476
$ws.Range("J15").Value2 = 1.040440084365198
$ws.Range("K15").Value2 = 1.046555097692649
$ws.Range("L15").Value2 = 1.041307176480238
$ws.Range("M15").Value2 = 1.054545724062241
$ws.Range("N15").Value2 = 1.041917628625142

$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.034603307661585
$ws.Range("D16").Value2 = 1.043591848192205
$ws.Range("E16").Value2 = 1.038438371296444
$ws.Range("F16").Value2 = 1.05169347550986
$ws.Range("I16").Value2 = 1.039274218258486
$ws.Range("J16").Value2 = 1.040929582793371
$ws.Range("K16").Value2 = 1.047002520611512
$ws.Range("L16").Value2 = 1.041867216106965
$ws.Range("M16").Value2 = 1.055075976004962
$ws.Range("N16").Value2 = 1.04240782219723

$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.035048439215466
$ws.Range("D17").Value2 = 1.043944464909154
$ws.Range("E17").Value2 = 1.038862522224616
$ws.Range("F17").Value2 = 1.05209768539758
$ws.Range("I17").Value2 = 1.039387393133757
$ws.Range("J17").Value2 = 1.041236438206313
$ws.Range("K17").Value2 = 1.04728285041639
$ws.Range("L17").Value2 = 1.042218416856808
$ws.Range("M17").Value2 = 1.055408368449536
$ws.Range("N17").Value2 = 1.042715113380055

$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.035308068059567
$ws.Range("D18").Value2 = 1.044150097337237
$ws.Range("E18").Value2 = 1.039109956858053
$ws.Range("F18").Value2 = 1.052333446133116
$ws.Range("I18").Value2 = 1.039453195784203
$ws.Range("J18").Value2 = 1.041415352041223
$ws.Range("K18").Value2 = 1.047446244401209
$ws.Range("L18").Value2 = 1.042423231446617
$ws.Range("M18").Value2 = 1.055602167970978
$ws.Range("N18").Value2 = 1.042894281293131

$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.035396593324674
$ws.Range("D19").Value2 = 1.044220205406219
$ws.Range("E19").Value2 = 1.039194331559707
$ws.Range("F19").Value2 = 1.052413833098303
$ws.Range("I19").Value2 = 1.039475597113913
$ws.Range("J19").Value2 = 1.041476345224506
$ws.Range("K19").Value2 = 1.047501937583435
$ws.Range("L19").Value2 = 1.042493062034308
$ws.Range("M19").Value2 = 1.055668235149618
$ws.Range("N19").Value2 = 1.042955361093731

$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.035000681738516
$ws.Range("D20").Value2 = 1.043906636915481
$ws.Range("E20").Value2 = 1.038817011273079
$ws.Range("F20").Value2 = 1.052054318337936
$ws.Range("I20").Value2 = 1.039375272297907
$ws.Range("J20").Value2 = 1.041203522741683
$ws.Range("K20").Value2 = 1.047252785852646
$ws.Range("L20").Value2 = 1.042180739950828
$ws.Range("M20").Value2 = 1.055372714098107
$ws.Range("N20").Value2 = 1.042682151171689

$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.033713427680339
$ws.Range("D21").Value2 = 1.042886690194556
$ws.Range("E21").Value2 = 1.037590711274143
$ws.Range("F21").Value2 = 1.050885405477197
$ws.Range("I21").Value2 = 1.039046634177077
$ws.Range("J21").Value2 = 1.040315729548209
$ws.Range("K21").Value2 = 1.046441385931506
$ws.Range("L21").Value2 = 1.0411649394451
$ws.Range("M21").Value2 = 1.054411012848775
$ws.Range("N21").Value2 = 1.04179309721005

$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.032903990727745
$ws.Range("D22").Value2 = 1.042245020569801
$ws.Range("E22").Value2 = 1.036819994319823
$ws.Range("F22").Value2 = 1.050150388054789
$ws.Range("I22").Value2 = 1.038838110565842
$ws.Range("J22").Value2 = 1.039756906494173
$ws.Range("K22").Value2 = 1.045930164209693
$ws.Range("L22").Value2 = 1.04052594686246
$ws.Range("M22").Value2 = 1.053805635287826
$ws.Range("N22").Value2 = 1.041233480563179

$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.033333097487055
$ws.Range("D23").Value2 = 1.042585218706962
$ws.Range("E23").Value2 = 1.037228537145142
$ws.Range("F23").Value2 = 1.050540042000098
$ws.Range("I23").Value2 = 1.038948832617161
$ws.Range("J23").Value2 = 1.040053209291477
$ws.Range("K23").Value2 = 1.046201273341715
$ws.Range("L23").Value2 = 1.040864719323103
$ws.Range("M23").Value2 = 1.054126625564379
$ws.Range("N23").Value2 = 1.041530204144441

$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.035022261313609
$ws.Range("D24").Value2 = 1.043923729889997
$ws.Range("E24").Value2 = 1.038837575605967
$ws.Range("F24").Value2 = 1.052073914070113
$ws.Range("I24").Value2 = 1.03938074983169
$ws.Range("J24").Value2 = 1.041218396038545
$ws.Range("K24").Value2 = 1.047266371098187
$ws.Range("L24").Value2 = 1.04219776463032
$ws.Range("M24").Value2 = 1.055388825009796
$ws.Range("N24").Value2 = 1.042697045590339

$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.036982099081089
$ws.Range("D25").Value2 = 1.045475314051032
$ws.Range("E25").Value2 = 1.040706133824381
$ws.Range("F25").Value2 = 1.053853572677482
$ws.Range("I25").Value2 = 1.039873714851434
$ws.Range("J25").Value2 = 1.042567796739132
$ws.Range("K25").Value2 = 1.048497743828878
$ws.Range("L25").Value2 = 1.043743323930611
$ws.Range("M25").Value2 = 1.05685042505555
$ws.Range("N25").Value2 = 1.044048362594699
